$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test")

# New headers for the added Properties File Logic columns
$ws.Range("L1").Value = "xLocation"
$ws.Range("M1").Value = "yLocation"

# Update FirstNameTextBox / SecondNameTextBox element-find-by-value entries
$ws.Range("D4").Value = "u_0_1"
$ws.Range("D5").Value = "u_0_3"

# New x/y location data for each row
$ws.Range("L2").Value = 872
$ws.Range("M2").Value = 23

$ws.Range("L3").Value = 1036
$ws.Range("M3").Value = 33

$ws.Range("L4").Value = 849
$ws.Range("M4").Value = 200

$ws.Range("L5").Value = 1054
$ws.Range("M5").Value = 200

# Match the selection shown in the diff (M5 now the last touched cell)
$ws.Range("M5").Select()
